$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted. In the published sheet this
# shows up as two brand-new data rows (114 and 115) with every subsequent
# row (old 114:137) pushed down by two positions (new 116:139).
$ws.Range("A114:A115").EntireRow.Insert()

# New row 114 - Packham's Triumph, Primera
$ws.Range("A114").Value = 4
$ws.Range("B114").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value = "Los Lagos"
$ws.Range("D114").Value = 44474
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100104
$ws.Range("H114").Value = "Frutos de pepita"
$ws.Range("I114").Value = 100104005
$ws.Range("J114").Value = "Pera"
$ws.Range("K114").Value = "Packham's Triumph"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 300
$ws.Range("N114").Value = 16000
$ws.Range("O114").Value = 16000
$ws.Range("P114").Value = 16000
$ws.Range("Q114").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 1067
$ws.Range("T114").Value = 15

# New row 115 - Packham's Triumph, Segunda
$ws.Range("A115").Value = 4
$ws.Range("B115").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C115").Value = "Los Lagos"
$ws.Range("D115").Value = 44474
$ws.Range("E115").Value = 10
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100104
$ws.Range("H115").Value = "Frutos de pepita"
$ws.Range("I115").Value = 100104005
$ws.Range("J115").Value = "Pera"
$ws.Range("K115").Value = "Packham's Triumph"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 200
$ws.Range("N115").Value = 13000
$ws.Range("O115").Value = 13000
$ws.Range("P115").Value = 13000
$ws.Range("Q115").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R115").Value = "Región de O'Higgins"
$ws.Range("S115").Value = 867
$ws.Range("T115").Value = 15
